# Insert a new price-report row at row 98 (pushing the existing row 98
# and everything below it down by one), then populate the new row with
# its data. This mirrors the source diff: dimension grows from A1:R178
# to A1:R179, and a brand-new record appears right before the former
# row 98 (which is now row 99).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 98..178 down to 99..179
$ws.Rows("98:98").Insert()

# Fill in the newly inserted row 98 with its values
$ws.Cells.Item(98, 1).Value = 11
$ws.Cells.Item(98, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(98, 3).Value = "Bíobío"
$ws.Cells.Item(98, 4).Value = 44944
$ws.Cells.Item(98, 5).Value = 8
$ws.Cells.Item(98, 6).Value = 100112043
$ws.Cells.Item(98, 7).Value = "Pepino ensalada"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 300
$ws.Cells.Item(98, 11).Value = 9000
$ws.Cells.Item(98, 12).Value = 9500
$ws.Cells.Item(98, 13).Value = 9250
$ws.Cells.Item(98, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(98, 15).Value = "Región Metropolitana"
$ws.Cells.Item(98, 16).Value = 154
$ws.Cells.Item(98, 17).Value = 60
$ws.Cells.Item(98, 18).Value = "Hortaliza"
